# Generate Report for Handoff
# Adds a new tracked file "a4494637-06dc-4c47-ab0c-923dbad5d284.md" as a row
# just above the trailing ".localization-config" row on all three sheets
# (Overview, zh-cn, de-de), matching the "Ready for handoff" rows already
# present for 4bd0e3d7-... and 878fd696-....

$wb = $excel.ActiveWorkbook

$missing = [System.Reflection.Missing]::Value

$newFileId = "a4494637-06dc-4c47-ab0c-923dbad5d284.md"
$newXlfHashZh = "a4494637-06dc-4c47-ab0c-923dbad5d284.87e8db1b0a0a9b4e1ab7e1a4f753d23b4291cd02.zh-cn.xlf"
$newXlfHashDe = "a4494637-06dc-4c47-ab0c-923dbad5d284.87e8db1b0a0a9b4e1ab7e1a4f753d23b4291cd02.de-de.xlf"
$newHandoffDtZh = "2016-02-24 07:39:23"
$newHandoffDtDe = "2016-02-24 07:39:36"

$mdCommit = "9732445c64dc794a76d759decee4ac1dcc9f6bf0"
$mdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/$mdCommit/e2e/$newFileId"
$cfgUrl = "https://github.com/OpenLocalizationTest/oltest/blob/$mdCommit/.localization-config"
$xlfZhUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/87e8db1b0a0a9b4e1ab7e1a4f753d23b4291cd02/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/$newXlfHashZh"
$xlfDeUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/87e8db1b0a0a9b4e1ab7e1a4f753d23b4291cd02/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/$newXlfHashDe"

# ---------------------------------------------------------------------------
# Sheet 1: "Overview"  (columns A=File Name, B=zh-cn status, C=de-de status)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")
$ws1.Rows.Item(7).Insert()
$ws1.Range("A7").Value = $newFileId
$ws1.Range("B7").Value = "Ready for handoff"
$ws1.Range("C7").Value = "Ready for handoff"

$ws1.Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/f8f0b1d0c1566b8db2aabfc939fca5fafe84e8ff/e2e/4a674e25-3ba2-4e52-833b-68918e322936.md", $missing, $missing, "4a674e25-3ba2-4e52-833b-68918e322936.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/041cd0bce2ae55136360588cd7a6585f43442cf6/e2e/6d2e9fb0-015e-48a4-991f-2a6b03475b8e.md", $missing, $missing, "6d2e9fb0-015e-48a4-991f-2a6b03475b8e.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/041cd0bce2ae55136360588cd7a6585f43442cf6/e2e/fe0528f8-4d62-44e2-a2a7-ea8e1baf428f.md", $missing, $missing, "fe0528f8-4d62-44e2-a2a7-ea8e1baf428f.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/9732445c64dc794a76d759decee4ac1dcc9f6bf0/e2e/4bd0e3d7-dd1b-4a5d-9d73-cd7a682f7367.md", $missing, $missing, "4bd0e3d7-dd1b-4a5d-9d73-cd7a682f7367.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/ba534de833b74f3c135e0c59685e16161fbb53d9/e2e/878fd696-c78d-4955-9f93-3473f0657199.md", $missing, $missing, "878fd696-c78d-4955-9f93-3473f0657199.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A7"), $mdUrl, $missing, $missing, $newFileId) | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A8"), $cfgUrl, $missing, $missing, ".localization-config") | Out-Null

# ---------------------------------------------------------------------------
# Sheet 2: "zh-cn"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws2.Rows.Item(7).Insert()
$ws2.Range("A7").Value = $newFileId
$ws2.Range("B7").Value = "Ready for handoff"
$ws2.Range("C7").Value = $newXlfHashZh
$ws2.Range("D7").Value = $newHandoffDtZh
$ws2.Range("G7").Value = "0001-01-01 00:00:00"
$ws2.Range("H7").Value = "Include"

$ws2.Range("D8").Value = "0001-01-01 00:00:00"

$ws2.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/f8f0b1d0c1566b8db2aabfc939fca5fafe84e8ff/e2e/4a674e25-3ba2-4e52-833b-68918e322936.md", $missing, $missing, "4a674e25-3ba2-4e52-833b-68918e322936.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/77985267f2f5d0f05cb15331da5b368c776df01f/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/4a674e25-3ba2-4e52-833b-68918e322936.f470f54e3993e463d5a8ee8c1e80c4044f0f6bd7.zh-cn.xlf", $missing, $missing, "4a674e25-3ba2-4e52-833b-68918e322936.f470f54e3993e463d5a8ee8c1e80c4044f0f6bd7.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/041cd0bce2ae55136360588cd7a6585f43442cf6/e2e/6d2e9fb0-015e-48a4-991f-2a6b03475b8e.md", $missing, $missing, "6d2e9fb0-015e-48a4-991f-2a6b03475b8e.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/7c63bc377c71a8525dddb5735dfbd23c36129465/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/6d2e9fb0-015e-48a4-991f-2a6b03475b8e.9b44a5e302e2c98dd79ec253ed6cb9040a68128e.zh-cn.xlf", $missing, $missing, "6d2e9fb0-015e-48a4-991f-2a6b03475b8e.9b44a5e302e2c98dd79ec253ed6cb9040a68128e.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/041cd0bce2ae55136360588cd7a6585f43442cf6/e2e/fe0528f8-4d62-44e2-a2a7-ea8e1baf428f.md", $missing, $missing, "fe0528f8-4d62-44e2-a2a7-ea8e1baf428f.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/7c63bc377c71a8525dddb5735dfbd23c36129465/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/fe0528f8-4d62-44e2-a2a7-ea8e1baf428f.70460a93df7465fe6297039e43cad0efe0da6720.zh-cn.xlf", $missing, $missing, "fe0528f8-4d62-44e2-a2a7-ea8e1baf428f.70460a93df7465fe6297039e43cad0efe0da6720.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/9732445c64dc794a76d759decee4ac1dcc9f6bf0/e2e/4bd0e3d7-dd1b-4a5d-9d73-cd7a682f7367.md", $missing, $missing, "4bd0e3d7-dd1b-4a5d-9d73-cd7a682f7367.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("C5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/7109907c301e84e01885eee76c99232259334523/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/4bd0e3d7-dd1b-4a5d-9d73-cd7a682f7367.57b603f9b8643304657877f8a4d27363813033dc.zh-cn.xlf", $missing, $missing, "4bd0e3d7-dd1b-4a5d-9d73-cd7a682f7367.57b603f9b8643304657877f8a4d27363813033dc.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/ba534de833b74f3c135e0c59685e16161fbb53d9/e2e/878fd696-c78d-4955-9f93-3473f0657199.md", $missing, $missing, "878fd696-c78d-4955-9f93-3473f0657199.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("C6"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/12afa8f78326ffaa12794d287ba02c21569722cc/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/878fd696-c78d-4955-9f93-3473f0657199.b770f6653a0056dfab75e6a68c30d6d72a99a4fc.zh-cn.xlf", $missing, $missing, "878fd696-c78d-4955-9f93-3473f0657199.b770f6653a0056dfab75e6a68c30d6d72a99a4fc.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A7"), $mdUrl, $missing, $missing, $newFileId) | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("C7"), $xlfZhUrl, $missing, $missing, $newXlfHashZh) | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A8"), $cfgUrl, $missing, $missing, ".localization-config") | Out-Null

# ---------------------------------------------------------------------------
# Sheet 3: "de-de"
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")
$ws3.Rows.Item(7).Insert()
$ws3.Range("A7").Value = $newFileId
$ws3.Range("B7").Value = "Ready for handoff"
$ws3.Range("C7").Value = $newXlfHashDe
$ws3.Range("D7").Value = $newHandoffDtDe
$ws3.Range("G7").Value = "0001-01-01 00:00:00"
$ws3.Range("H7").Value = "Include"

$ws3.Range("D8").Value = "0001-01-01 00:00:00"

$ws3.Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/f8f0b1d0c1566b8db2aabfc939fca5fafe84e8ff/e2e/4a674e25-3ba2-4e52-833b-68918e322936.md", $missing, $missing, "4a674e25-3ba2-4e52-833b-68918e322936.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a0a97eed87c644a32aa7a5bb88df45f2b96f2d48/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/4a674e25-3ba2-4e52-833b-68918e322936.f470f54e3993e463d5a8ee8c1e80c4044f0f6bd7.de-de.xlf", $missing, $missing, "4a674e25-3ba2-4e52-833b-68918e322936.f470f54e3993e463d5a8ee8c1e80c4044f0f6bd7.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/041cd0bce2ae55136360588cd7a6585f43442cf6/e2e/6d2e9fb0-015e-48a4-991f-2a6b03475b8e.md", $missing, $missing, "6d2e9fb0-015e-48a4-991f-2a6b03475b8e.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/dfd61c89e4eb849ed0b082c6f78e03cb463c6238/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/6d2e9fb0-015e-48a4-991f-2a6b03475b8e.9b44a5e302e2c98dd79ec253ed6cb9040a68128e.de-de.xlf", $missing, $missing, "6d2e9fb0-015e-48a4-991f-2a6b03475b8e.9b44a5e302e2c98dd79ec253ed6cb9040a68128e.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/041cd0bce2ae55136360588cd7a6585f43442cf6/e2e/fe0528f8-4d62-44e2-a2a7-ea8e1baf428f.md", $missing, $missing, "fe0528f8-4d62-44e2-a2a7-ea8e1baf428f.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/dfd61c89e4eb849ed0b082c6f78e03cb463c6238/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/fe0528f8-4d62-44e2-a2a7-ea8e1baf428f.70460a93df7465fe6297039e43cad0efe0da6720.de-de.xlf", $missing, $missing, "fe0528f8-4d62-44e2-a2a7-ea8e1baf428f.70460a93df7465fe6297039e43cad0efe0da6720.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/9732445c64dc794a76d759decee4ac1dcc9f6bf0/e2e/4bd0e3d7-dd1b-4a5d-9d73-cd7a682f7367.md", $missing, $missing, "4bd0e3d7-dd1b-4a5d-9d73-cd7a682f7367.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("C5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/40938d89938ed1a24dda37daf26d4dad67d1621c/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/4bd0e3d7-dd1b-4a5d-9d73-cd7a682f7367.57b603f9b8643304657877f8a4d27363813033dc.de-de.xlf", $missing, $missing, "4bd0e3d7-dd1b-4a5d-9d73-cd7a682f7367.57b603f9b8643304657877f8a4d27363813033dc.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/ba534de833b74f3c135e0c59685e16161fbb53d9/e2e/878fd696-c78d-4955-9f93-3473f0657199.md", $missing, $missing, "878fd696-c78d-4955-9f93-3473f0657199.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("C6"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/7a65ed68eea20202716b2037f6517358137edb75/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/878fd696-c78d-4955-9f93-3473f0657199.b770f6653a0056dfab75e6a68c30d6d72a99a4fc.de-de.xlf", $missing, $missing, "878fd696-c78d-4955-9f93-3473f0657199.b770f6653a0056dfab75e6a68c30d6d72a99a4fc.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A7"), $mdUrl, $missing, $missing, $newFileId) | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("C7"), $xlfDeUrl, $missing, $missing, $newXlfHashDe) | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A8"), $cfgUrl, $missing, $missing, ".localization-config") | Out-Null

Write-Host "Done generating handoff report rows."
